$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 6271.3313874533
$ws.Range("C4").Value = 31.14975803354112
$ws.Range("D6").Value = 9474.236236812256
